$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("MapNpcInfoData")
$ws2 = $wb.Worksheets.Item("MapNpcMenuData")

# Update the id-column header strings (NpcId key column type annotation)
$ws2.Range("A1").Value = "int;id"
$ws1.Range("A1").Value = "int;id;key"

# Update selections and active sheet/tab to match the new view state
$ws2.Range("A8").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("B6").Select() | Out-Null
